# Correcting some documentation issues
#
# 1) Rename the worksheet tab to match the workbook's file name (H21R00 -> H3BR20)
# 2) Clear the highlight fill that was put on the "D1" designator row (A20) -
#    it used the Accent5 theme tint; it should go back to the plain
#    Background 1 (white/no-highlight) theme color.
# 3) Update the saved cursor/selection position (scrolled back to the top,
#    with the active cell now on B17 instead of A15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the sheet name so it matches the board/file name.
$ws.Name = "H3BR20"

# 2) Remove the review-highlight color from A20 (the D1 / display row) -
#    set its interior back to the theme's Background 1 (white).
$ws.Range("A20").Interior.ThemeColor = 2

# 3) Restore the view: scroll back to the top of the sheet and leave the
#    selection on B17.
$ws.Range("B17").Select()
